$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The header row labels (A1, C1, E1) identified the row as test-block
# messages ("test_msg_xx"). They are being renamed to identify block
# messages ("block_msg_xx") instead.
$ws.Range("A1").Value = "block_msg_en"
$ws.Range("C1").Value = "block_msg_es"
$ws.Range("E1").Value = "block_msg_fr"
